$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.353155374526978
$ws.Range("B1").Value = 2.677488565444946
$ws.Range("C1").Value = 4.175845146179199
$ws.Range("D1").Value = 4.245255947113037
$ws.Range("E1").Value = 1.70209801197052
